# Ingreso - cargado.xlsx : add two new test-case rows (124, 125) to Hoja1
# and move the sheet view/selection down to the new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New shared strings must be introduced in the same order they appear in the
# target workbook: fechahora(124), fechahora(125), motivo(124), motivo(125).
$ws.Range("A124").Value = "2017-09-25-20:09:58"
$ws.Range("A125").Value = "2017-10-30-13:22:46"
$ws.Range("J124").Value = "ESTA-DOCU"
$ws.Range("J125").Value = "ALCO-DOCU"

# ---------------------------------------------------------------------------
# Row 124 (remaining columns)
# ---------------------------------------------------------------------------
$ws.Range("B124").Value = 3
$ws.Range("C124").Value = 5
$ws.Range("D124").Value = "S/P"
$ws.Range("E124").Value = "CHEVROLET"
$ws.Range("F124").Value = "CHEVY"
$ws.Range("F124").HorizontalAlignment = -4131
$ws.Range("G124").Value = "Negro"
$ws.Range("K124").Value = "Arenales"
$ws.Range("L124").Value = "S/D"
$ws.Range("M124").Value = "Venezuela"
$ws.Range("N124").Value = "Carlos Calvo"
$ws.Range("O124").Value = 8540
$ws.Range("P124").Value = 21
$ws.Range("Q124").Value = 39
$ws.Range("R124").Value = 25
$ws.Range("S124").Value = "S/D"
$ws.Range("T124").Value = 1
$ws.Range("U124").Value = 16855456
$ws.Range("U124").HorizontalAlignment = -4131
$ws.Range("AA124").Value = 16017888
$ws.Range("AB124").Value = 16085667
$ws.Range("AC124").Value = 565699
$ws.Range("AD124").Value = 680669
# Copy (rather than re-type) the literal text "true" so it stays a shared
# string instead of Excel auto-coercing the literal into a real Boolean.
$ws.Range("AE114").Copy($ws.Range("AE124"))
$ws.Range("AF124").Value = 1
$ws.Range("AG124").Value = 1
$ws.Range("AE114").Copy($ws.Range("AH124"))

# ---------------------------------------------------------------------------
# Row 125 (remaining columns)
# ---------------------------------------------------------------------------
$ws.Range("B125").Value = 3
$ws.Range("C125").Value = 1
$ws.Range("D125").Value = "B1995958"
$ws.Range("E125").Value = "FORD"
$ws.Range("F125").Value = "SIERRA"
$ws.Range("F125").HorizontalAlignment = -4131
$ws.Range("G125").Value = "Azul"
$ws.Range("K125").Value = "Ecuador"
$ws.Range("L125").Value = "S/D"
$ws.Range("M125").Value = "Saavedra"
$ws.Range("N125").Value = "Solis"
$ws.Range("O125").Value = 3606
$ws.Range("P125").Value = 22
$ws.Range("Q125").Value = 40
$ws.Range("R125").Value = 146
$ws.Range("S125").Value = "S/D"
$ws.Range("T125").Value = 1
$ws.Range("U125").Value = 16855456
$ws.Range("U125").HorizontalAlignment = -4131
$ws.Range("AA125").Value = 16017888
$ws.Range("AB125").Value = 16085667
$ws.Range("AC125").Value = 565699
$ws.Range("AD125").Value = 680669
$ws.Range("AE114").Copy($ws.Range("AE125"))
$ws.Range("AF125").Value = 1
$ws.Range("AG125").Value = 1
$ws.Range("AE114").Copy($ws.Range("AH125"))

# ---------------------------------------------------------------------------
# View: scroll down to the new rows and select the last one entered (J125)
# ---------------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 95
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("J125").Select()
